$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 1267
$ws1.Range("F6").Value = 17783
$ws1.Range("F13").Value = 9
$ws1.Range("F19").Value = 183
$ws1.Range("F25").Value = 260
$ws1.Range("F26").Value = 965
$ws1.Range("F31").Value = 11905
$ws1.Range("F34").Value = 195

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 1267
$ws4.Range("F6").Value = 17783
$ws4.Range("F13").Value = 9
$ws4.Range("F19").Value = 183
$ws4.Range("F25").Value = 260
$ws4.Range("F26").Value = 965
$ws4.Range("F33").Value = 11905
$ws4.Range("F36").Value = 195
